# PlayerPerformance_4325.xlsx edit
# - add "Player Info" sheet (before "ODI Batting")
# - add "ODI Batting Extra" sheet (after "ODI Batting")
# - on "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, convert the
#   scorecard URL values down to just the numeric match code, and drop
#   the handful of blank INNING_NUMBER cells (rows where the player did
#   not bat) entirely rather than leaving an empty string cell.

$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------
# 1. Create the two extra sheets, in the right tab order.
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

# re-resolve "ODI Batting" - inserting a sheet can shift positions
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiExtra = $wb.Worksheets.Add($null, $odiBatting)
$odiExtra.Name = "ODI Batting Extra"

# match the page margins used by the rest of the workbook (0.75in sides,
# 1in top/bottom, 0.5in header/footer == 54/54/72/72/36/36 points)
foreach ($sheet in @($playerInfo, $odiExtra)) {
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# ---------------------------------------------------------------------
# 2. "Player Info" sheet content.
# ---------------------------------------------------------------------
$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $piHeaders.Length; $col++) {
    $cell = $playerInfo.Cells.Item(1, $col)
    $cell.Value = $piHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4325"
$playerInfo.Cells.Item(2, 2).Value = "Jason Jonathan Roy"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ---------------------------------------------------------------------
# 3. "ODI Batting" sheet updates.
# ---------------------------------------------------------------------
$odiBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRow = $odiBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $odiBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    $idx = $url.IndexOf("MatchCode=")
    if ($idx -ge 0) {
        $code = $url.Substring($idx + 10)
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# Rows where INNING_NUMBER (col B) was only a placeholder empty string -
# drop the cell entirely (player did not bat in that match).
$blankInningRows = @(2, 74, 81, 87, 99, 108)
foreach ($r in $blankInningRows) {
    $odiBatting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------
# 4. "ODI Batting Extra" sheet content.
# ---------------------------------------------------------------------
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $cell = $odiExtra.Cells.Item(1, $col)
    $cell.Value = $extraHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4470", "1", "10", "0", "24.59%", "NO"),
    @("4471", "",  "",   "",  "",       "NO"),
    @("4598", "",  "",   "",  "",       "NO"),
    @("4599", "",  "",   "",  "",       "NO"),
    @("4602", "1", "15", "0", "40.73%", "YES"),
    @("4609", "1", "0",  "0", "",       "NO"),
    @("4613", "1", "2",  "1", "9.35%",  "NO"),
    @("4618", "1", "7",  "0", "15.83%", "NO"),
    @("4619", "",  "",   "",  "",       "NO"),
    @("4620", "1", "3",  "0", "6.97%",  "NO"),
    @("4622", "",  "",   "",  "",       "NO"),
    @("4660", "",  "",   "",  "",       "NO"),
    @("4663", "",  "",   "",  "",       "NO"),
    @("4666", "",  "",   "",  "",       "NO"),
    @("4698", "1", "11", "4", "41.70%", "NO"),
    @("4699", "1", "1",  "0", "2.63%",  "NO"),
    @("4700", "1", "0",  "0", "0.29%",  "NO"),
    @("4711", "1", "1",  "0", "1.89%",  "NO"),
    @("4713", "1", "18", "1", "40.49%", "YES"),
    @("4717", "1", "3",  "0", "9.69%",  "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $r = $i + 2
    $row = $extraRows[$i]

    $codeCell = $odiExtra.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    if ($row[1] -ne "") {
        $odiExtra.Cells.Item($r, 2).Value = [int]$row[1]
    }
    if ($row[2] -ne "") {
        $c = $odiExtra.Cells.Item($r, 3)
        $c.NumberFormat = "@"
        $c.Value = $row[2]
    }
    if ($row[3] -ne "") {
        $c = $odiExtra.Cells.Item($r, 4)
        $c.NumberFormat = "@"
        $c.Value = $row[3]
    }
    if ($row[4] -ne "") {
        $c = $odiExtra.Cells.Item($r, 5)
        $c.NumberFormat = "@"
        $c.Value = $row[4]
    }
    $odiExtra.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------------
# 5. Leave "Player Info" as the active sheet/tab (matches activeTab="0").
# ---------------------------------------------------------------------
$playerInfo.Activate()
[void]$playerInfo.Range("A1").Select()
